$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 18 (shift rows up) since new data only has 16 data rows (rows 2-17)
$ws.Rows.Item(18).Delete() | Out-Null

# Row 2 (column A untouched - unchanged by diff)
$ws.Range("B2").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C2").Value2 = "HEY1"
$ws.Range("D2").Value2 = "2/12"
$ws.Range("E2").Value2 = 0.000044656194172606244
$ws.Range("F2").Value2 = 0.012816327727537992
$ws.Range("G2").Value2 = 196.07843137254903
$ws.Range("H2").Value2 = 1964.0230456849977
$ws.Range("I2").Value2 = "SMAD3;FBXW7"
$ws.Range("J2").Value2 = 2.0

# Row 3 (column A untouched - unchanged by diff)
$ws.Range("B3").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C3").Value2 = "KLF5"
$ws.Range("D3").Value2 = "2/35"
$ws.Range("E3").Value2 = 0.0003979873115628295
$ws.Range("F3").Value2 = 0.038074119472844024
$ws.Range("G3").Value2 = 67.22689075630252
$ws.Range("H3").Value2 = 526.3254073064493
$ws.Range("I3").Value2 = "SMAD3;FBXW7"
$ws.Range("J3").Value2 = 2.0

# Row 4 (column A untouched - unchanged by diff)
$ws.Range("B4").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C4").Value2 = "KLF4"
$ws.Range("D4").Value2 = "2/31"
$ws.Range("E4").Value2 = 0.0003116539576784677
$ws.Range("F4").Value2 = 0.044722342926860115
$ws.Range("G4").Value2 = 75.90132827324479
$ws.Range("H4").Value2 = 612.7982615314888
$ws.Range("I4").Value2 = "SMAD3;APC"
$ws.Range("J4").Value2 = 2.0

# Row 5 (column A untouched - unchanged by diff)
$ws.Range("B5").Value2 = "WikiPathways_2019_Human"
$ws.Range("C5").Value2 = "Factors and pathways affecting insulin-like growth factor (IGF1)-Akt signaling WP3850"
$ws.Range("D5").Value2 = "2/31"
$ws.Range("E5").Value2 = 0.0003116539576784677
$ws.Range("F5").Value2 = 0.04903355600807891
$ws.Range("G5").Value2 = 75.90132827324479
$ws.Range("H5").Value2 = 612.7982615314888
$ws.Range("I5").Value2 = "SMAD3;NEB"
$ws.Range("J5").Value2 = 2.0

# Row 6 (column A untouched - unchanged by diff)
$ws.Range("B6").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C6").Value2 = "JUN"
$ws.Range("D6").Value2 = "3/248"
$ws.Range("E6").Value2 = 0.001126273578445804
$ws.Range("F6").Value2 = 0.06464810340278915
$ws.Range("G6").Value2 = 14.231499051233397
$ws.Range("H6").Value2 = 96.61538160204768
$ws.Range("I6").Value2 = "SMAD3;APC;FBXW7"
$ws.Range("J6").Value2 = 3.0

# Row 7 (column A untouched - unchanged by diff)
$ws.Range("B7").Value2 = "WikiPathways_2019_Human"
$ws.Range("C7").Value2 = "Extracellular vesicle-mediated signaling in recipient cells WP2870"
$ws.Range("D7").Value2 = "2/30"
$ws.Range("E7").Value2 = 0.00029169286205301905
$ws.Range("F7").Value2 = 0.0688395154445125
$ws.Range("G7").Value2 = 78.43137254901961
$ws.Range("H7").Value2 = 638.416404043489
$ws.Range("I7").Value2 = "SMAD3;APC"
$ws.Range("J7").Value2 = 2.0

# Row 8 (column A untouched - unchanged by diff)
$ws.Range("B8").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C8").Value2 = "SMAD4"
$ws.Range("D8").Value2 = "3/245"
$ws.Range("E8").Value2 = 0.0010874447654406513
$ws.Range("F8").Value2 = 0.07802416192036674
$ws.Range("G8").Value2 = 14.40576230492197
$ws.Range("H8").Value2 = 98.30383558257928
$ws.Range("I8").Value2 = "SMAD3;APC;THOC2"
$ws.Range("J8").Value2 = 3.0

# Row 9 (column A untouched - unchanged by diff)
$ws.Range("B9").Value2 = "WikiPathways_2019_Human"
$ws.Range("C9").Value2 = "Hypothesized Pathways in Pathogenesis of Cardiovascular Disease WP3668"
$ws.Range("D9").Value2 = "2/25"
$ws.Range("E9").Value2 = 0.00020167031602328324
$ws.Range("F9").Value2 = 0.09518838916298969
$ws.Range("G9").Value2 = 94.11764705882354
$ws.Range("H9").Value2 = 800.8354157912519
$ws.Range("I9").Value2 = "FBN2;SMAD3"
$ws.Range("J9").Value2 = 2.0

# Row 10 (column A untouched - unchanged by diff)
$ws.Range("B10").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C10").Value2 = "TP63"
$ws.Range("D10").Value2 = "2/120"
$ws.Range("E10").Value2 = 0.004577501584745696
$ws.Range("F10").Value2 = 0.11943117771109225
$ws.Range("G10").Value2 = 19.607843137254903
$ws.Range("H10").Value2 = 105.61964578847743
$ws.Range("I10").Value2 = "SMAD3;FBXW7"
$ws.Range("J10").Value2 = 2.0

# Row 11 (column A untouched - unchanged by diff)
$ws.Range("B11").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C11").Value2 = "HIF1A"
$ws.Range("D11").Value2 = "2/126"
$ws.Range("E11").Value2 = 0.0050336357625974625
$ws.Range("F11").Value2 = 0.12038778865545598
$ws.Range("G11").Value2 = 18.674136321195146
$ws.Range("H11").Value2 = 98.81629767227176
$ws.Range("I11").Value2 = "SMAD3;FBXW7"
$ws.Range("J11").Value2 = 2.0

# Row 12 (column A untouched - unchanged by diff)
$ws.Range("B12").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C12").Value2 = "CTNNB1"
$ws.Range("D12").Value2 = "3/398"
$ws.Range("E12").Value2 = 0.004321924909294262
$ws.Range("F12").Value2 = 0.12403924489674531
$ws.Range("G12").Value2 = 8.867868755542418
$ws.Range("H12").Value2 = 48.277159873892714
$ws.Range("I12").Value2 = "SMAD3;APC;LRRK2"
$ws.Range("J12").Value2 = 3.0

# Row 13 (column A untouched - unchanged by diff)
$ws.Range("B13").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C13").Value2 = "CEBPA"
$ws.Range("D13").Value2 = "2/113"
$ws.Range("E13").Value2 = 0.004071107650863602
$ws.Range("F13").Value2 = 0.12982309953309484
$ws.Range("G13").Value2 = 20.82248828735034
$ws.Range("H13").Value2 = 114.6036474016536
$ws.Range("I13").Value2 = "SMAD3;FBXW7"
$ws.Range("J13").Value2 = 2.0

# Row 14 (column A untouched - unchanged by diff)
$ws.Range("B14").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C14").Value2 = "VDR"
$ws.Range("D14").Value2 = "2/109"
$ws.Range("E14").Value2 = 0.0037943088288471378
$ws.Range("F14").Value2 = 0.13612082923489108
$ws.Range("G14").Value2 = 21.586616297895308
$ws.Range("H14").Value2 = 120.3292609068028
$ws.Range("I14").Value2 = "SMAD3;CFH"
$ws.Range("J14").Value2 = 2.0

# Row 15 (column A untouched - unchanged by diff)
$ws.Range("B15").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C15").Value2 = "EPAS1"
$ws.Range("D15").Value2 = "2/106"
$ws.Range("E15").Value2 = 0.003592753073271419
$ws.Range("F15").Value2 = 0.14730287600412817
$ws.Range("G15").Value2 = 22.197558268590456
$ws.Range("H15").Value2 = 124.94642613986858
$ws.Range("I15").Value2 = "SMAD3;APC"
$ws.Range("J15").Value2 = 2.0

# Row 16 (column A untouched - unchanged by diff)
$ws.Range("B16").Value2 = "Reactome_2016"
$ws.Range("C16").Value2 = "O-linked glycosylation_Homo sapiens_R-HSA-5173105"
$ws.Range("D16").Value2 = "3/110"
$ws.Range("E16").Value2 = 0.00010406304392893897
$ws.Range("F16").Value2 = 0.15921645721127664
$ws.Range("G16").Value2 = 32.0855614973262
$ws.Range("H16").Value2 = 294.241079706687
$ws.Range("I16").Value2 = "WBSCR17;MUC16;ADAMTS9"
$ws.Range("J16").Value2 = 3.0

# Row 17 (column A untouched - unchanged by diff)
$ws.Range("B17").Value2 = "Transcription_Factor_PPIs"
$ws.Range("C17").Value2 = "ESR2"
$ws.Range("D17").Value2 = "3/365"
$ws.Range("E17").Value2 = 0.003389579926121341
$ws.Range("F17").Value2 = 0.16213490646613746
$ws.Range("G17").Value2 = 9.6696212731668
$ws.Range("H17").Value2 = 54.99161270781345
$ws.Range("I17").Value2 = "FLG;SMAD3;APC"
$ws.Range("J17").Value2 = 3.0
